$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value while preventing Excel's automatic
# date/time/number inference (e.g. "2023-09-08", "19:29"), without leaving
# any residual number-format/style on the cell.
function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 25
$ws.Range("A25").Value = 112013691
$ws.Range("B25").Value = 88489
$ws.Range("C25").Value = "Ovaliderad"
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 1962
$ws.Range("F25").Value = "Vaddporing"
$ws.Range("G25").Value = "Anomoporia kamtschatica"
$ws.Range("H25").Value = "(Parmasto) Bondartseva"
$ws.Range("P25").Value = "Spångmyran, Ås lm"
$ws.Range("Q25").Value = 610134.4051595986
$ws.Range("R25").Value = 7121460.896015909
$ws.Range("S25").Value = 25
$ws.Range("T25").Value = "Västerbotten"
$ws.Range("U25").Value = "Åsele"
$ws.Range("V25").Value = "Åsele lappmark"
$ws.Range("W25").Value = "Åsele"
Set-TextCell "Y25" "2023-09-08"
Set-TextCell "Z25" "19:29"
Set-TextCell "AA25" "2023-09-08"
Set-TextCell "AB25" "19:29"
$ws.Range("AD25").Value = $false
$ws.Range("AE25").Value = $false
$ws.Range("AG25").Value = $false
$ws.Range("AW25").Value = "Isak Vahlström"
$ws.Range("AX25").Value = "Isak Vahlström"

# Row 26
$ws.Range("A26").Value = 112013700
$ws.Range("B26").Value = 77515
$ws.Range("C26").Value = "Ovaliderad"
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 6425
$ws.Range("F26").Value = "Garnlav"
$ws.Range("G26").Value = "Alectoria sarmentosa"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("P26").Value = "Spångmyran, Ås lm"
$ws.Range("Q26").Value = 610101.9650201321
$ws.Range("R26").Value = 7121415.702941997
$ws.Range("S26").Value = 25
$ws.Range("T26").Value = "Västerbotten"
$ws.Range("U26").Value = "Åsele"
$ws.Range("V26").Value = "Åsele lappmark"
$ws.Range("W26").Value = "Åsele"
Set-TextCell "Y26" "2023-09-08"
Set-TextCell "Z26" "19:35"
Set-TextCell "AA26" "2023-09-08"
Set-TextCell "AB26" "19:35"
$ws.Range("AD26").Value = $false
$ws.Range("AE26").Value = $false
$ws.Range("AG26").Value = $false
$ws.Range("AW26").Value = "Isak Vahlström"
$ws.Range("AX26").Value = "Isak Vahlström"

# Row 27
$ws.Range("A27").Value = 112013697
$ws.Range("B27").Value = 89423
$ws.Range("C27").Value = "Ovaliderad"
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 5432
$ws.Range("F27").Value = "Granticka"
$ws.Range("G27").Value = "Porodaedalea chrysoloma"
$ws.Range("H27").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("P27").Value = "Spångmyran, Ås lm"
$ws.Range("Q27").Value = 610102.0736959254
$ws.Range("R27").Value = 7121412.654772604
$ws.Range("S27").Value = 25
$ws.Range("T27").Value = "Västerbotten"
$ws.Range("U27").Value = "Åsele"
$ws.Range("V27").Value = "Åsele lappmark"
$ws.Range("W27").Value = "Åsele"
Set-TextCell "Y27" "2023-09-08"
Set-TextCell "Z27" "19:35"
Set-TextCell "AA27" "2023-09-08"
Set-TextCell "AB27" "19:35"
$ws.Range("AD27").Value = $false
$ws.Range("AE27").Value = $false
$ws.Range("AG27").Value = $false
$ws.Range("AW27").Value = "Isak Vahlström"
$ws.Range("AX27").Value = "Isak Vahlström"

# Row 28
$ws.Range("A28").Value = 112013690
$ws.Range("B28").Value = 88489
$ws.Range("C28").Value = "Ovaliderad"
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 1962
$ws.Range("F28").Value = "Vaddporing"
$ws.Range("G28").Value = "Anomoporia kamtschatica"
$ws.Range("H28").Value = "(Parmasto) Bondartseva"
$ws.Range("P28").Value = "Spångmyran, Ås lm"
$ws.Range("Q28").Value = 610051.8565798617
$ws.Range("R28").Value = 7121425.252971379
$ws.Range("S28").Value = 25
$ws.Range("T28").Value = "Västerbotten"
$ws.Range("U28").Value = "Åsele"
$ws.Range("V28").Value = "Åsele lappmark"
$ws.Range("W28").Value = "Åsele"
Set-TextCell "Y28" "2023-09-08"
Set-TextCell "Z28" "19:43"
Set-TextCell "AA28" "2023-09-08"
Set-TextCell "AB28" "19:43"
$ws.Range("AD28").Value = $false
$ws.Range("AE28").Value = $false
$ws.Range("AG28").Value = $false
$ws.Range("AW28").Value = "Isak Vahlström"
$ws.Range("AX28").Value = "Isak Vahlström"

# Row 29
$ws.Range("A29").Value = 112013704
$ws.Range("B29").Value = 81248
$ws.Range("C29").Value = "Ovaliderad"
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 1312
$ws.Range("F29").Value = "Gammelgransskål"
$ws.Range("G29").Value = "Pseudographis pinicola"
$ws.Range("H29").Value = "(Nyl.) Rehm"
$ws.Range("P29").Value = "Spångmyran, Ås lm"
$ws.Range("Q29").Value = 610093.591720929
$ws.Range("R29").Value = 7121454.644715369
$ws.Range("S29").Value = 25
$ws.Range("T29").Value = "Västerbotten"
$ws.Range("U29").Value = "Åsele"
$ws.Range("V29").Value = "Åsele lappmark"
$ws.Range("W29").Value = "Åsele"
Set-TextCell "Y29" "2023-09-08"
Set-TextCell "Z29" "19:49"
Set-TextCell "AA29" "2023-09-08"
Set-TextCell "AB29" "19:49"
$ws.Range("AD29").Value = $false
$ws.Range("AE29").Value = $false
$ws.Range("AG29").Value = $false
$ws.Range("AW29").Value = "Isak Vahlström"
$ws.Range("AX29").Value = "Isak Vahlström"

# Row 30
$ws.Range("A30").Value = 112013696
$ws.Range("B30").Value = 86961
$ws.Range("C30").Value = "Ovaliderad"
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 4962
$ws.Range("F30").Value = "Mjölsvärting"
$ws.Range("G30").Value = "Lyophyllum semitale"
$ws.Range("H30").Value = "(Fr. : Fr.) Kühner"
$ws.Range("P30").Value = "Spångmyran, Ås lm"
$ws.Range("Q30").Value = 610070.1349689787
$ws.Range("R30").Value = 7121402.360087069
$ws.Range("S30").Value = 25
$ws.Range("T30").Value = "Västerbotten"
$ws.Range("U30").Value = "Åsele"
$ws.Range("V30").Value = "Åsele lappmark"
$ws.Range("W30").Value = "Åsele"
Set-TextCell "Y30" "2023-09-08"
Set-TextCell "Z30" "19:40"
Set-TextCell "AA30" "2023-09-08"
Set-TextCell "AB30" "19:40"
$ws.Range("AD30").Value = $false
$ws.Range("AE30").Value = $false
$ws.Range("AG30").Value = $false
$ws.Range("AW30").Value = "Isak Vahlström"
$ws.Range("AX30").Value = "Isak Vahlström"

# Row 31
$ws.Range("A31").Value = 112013698
$ws.Range("B31").Value = 77515
$ws.Range("C31").Value = "Ovaliderad"
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 6425
$ws.Range("F31").Value = "Garnlav"
$ws.Range("G31").Value = "Alectoria sarmentosa"
$ws.Range("H31").Value = "(Ach.) Ach."
$ws.Range("P31").Value = "Spångmyran, Ås lm"
$ws.Range("Q31").Value = 610094.4326785516
$ws.Range("R31").Value = 7121455.546697079
$ws.Range("S31").Value = 25
$ws.Range("T31").Value = "Västerbotten"
$ws.Range("U31").Value = "Åsele"
$ws.Range("V31").Value = "Åsele lappmark"
$ws.Range("W31").Value = "Åsele"
Set-TextCell "Y31" "2023-09-08"
Set-TextCell "Z31" "19:49"
Set-TextCell "AA31" "2023-09-08"
Set-TextCell "AB31" "19:49"
$ws.Range("AD31").Value = $false
$ws.Range("AE31").Value = $false
$ws.Range("AG31").Value = $false
$ws.Range("AW31").Value = "Isak Vahlström"
$ws.Range("AX31").Value = "Isak Vahlström"

# Row 32
$ws.Range("A32").Value = 112013699
$ws.Range("B32").Value = 77515
$ws.Range("C32").Value = "Ovaliderad"
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 6425
$ws.Range("F32").Value = "Garnlav"
$ws.Range("G32").Value = "Alectoria sarmentosa"
$ws.Range("H32").Value = "(Ach.) Ach."
$ws.Range("P32").Value = "Spångmyran, Ås lm"
$ws.Range("Q32").Value = 610068.1736430819
$ws.Range("R32").Value = 7121408.394281525
$ws.Range("S32").Value = 25
$ws.Range("T32").Value = "Västerbotten"
$ws.Range("U32").Value = "Åsele"
$ws.Range("V32").Value = "Åsele lappmark"
$ws.Range("W32").Value = "Åsele"
Set-TextCell "Y32" "2023-09-08"
Set-TextCell "Z32" "19:40"
Set-TextCell "AA32" "2023-09-08"
Set-TextCell "AB32" "19:40"
$ws.Range("AD32").Value = $false
$ws.Range("AE32").Value = $false
$ws.Range("AG32").Value = $false
$ws.Range("AW32").Value = "Isak Vahlström"
$ws.Range("AX32").Value = "Isak Vahlström"

# Row 33
$ws.Range("A33").Value = 112013703
$ws.Range("B33").Value = 77515
$ws.Range("C33").Value = "Ovaliderad"
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 6425
$ws.Range("F33").Value = "Garnlav"
$ws.Range("G33").Value = "Alectoria sarmentosa"
$ws.Range("H33").Value = "(Ach.) Ach."
$ws.Range("P33").Value = "Spångmyran, Ås lm"
$ws.Range("Q33").Value = 610144.4332068264
$ws.Range("R33").Value = 7121461.253672058
$ws.Range("S33").Value = 25
$ws.Range("T33").Value = "Västerbotten"
$ws.Range("U33").Value = "Åsele"
$ws.Range("V33").Value = "Åsele lappmark"
$ws.Range("W33").Value = "Åsele"
Set-TextCell "Y33" "2023-09-08"
Set-TextCell "Z33" "19:28"
Set-TextCell "AA33" "2023-09-08"
Set-TextCell "AB33" "19:28"
$ws.Range("AD33").Value = $false
$ws.Range("AE33").Value = $false
$ws.Range("AG33").Value = $false
$ws.Range("AW33").Value = "Isak Vahlström"
$ws.Range("AX33").Value = "Isak Vahlström"

